# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") holds a per-row strike count that is recomputed from the
# underlying TB/PC/dS0/dSF data. This applies the regenerated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 3
    5  = 1
    6  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 2
    26 = 2
    29 = 1
    30 = 2
    32 = 1
    33 = 3
    34 = 1
    35 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
